$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoginData")

# Row 4: test3@yahoo.com / pass2322 / N
$ws.Range("A4").Value = "test3@yahoo.com"
$ws.Range("B4").Value = "pass2322"

# Row 5: test4@hotmail.com / pass121 / Y
$ws.Range("A5").Value = "test4@hotmail.com"
$ws.Range("B5").Value = "pass121"

# Runmode column (shared string order matches the target sst table)
$ws.Range("C4").Value = "N"
$ws.Range("C5").Value = "Y"

[void]$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:test3@yahoo.com")
[void]$ws.Hyperlinks.Add($ws.Range("A5"), "mailto:test4@hotmail.com")

# Match the existing "Username" hyperlink cell formatting (style index 2)
# by copying the format from A2 instead of leaving Excel's freshly minted style.
$ws.Range("A2").Copy()
[void]$ws.Range("A4").PasteSpecial(-4122)
[void]$ws.Range("A5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

[void]$ws.Range("B10").Select()
